{"js": "// Replace the date and each \"three-digit \u00f7 one-digit\" practice answer\n// with the updated values, as described in the commit diff.\n// Each old value is unique in the document, so a simple exact search\n// and replace (matchCase + wholeWords off is fine since strings are unique)\n// is sufficient and keeps existing run formatting (font/size) intact.\n\nconst replacements = [\n  [\"2024-02-13 Tuesday\", \"2024-02-14 Wednesday\"],\n  [\"160\u00f74=40, 0\", \"979\u00f75=195, 4\"],\n  [\"629\u00f78=78, 5\", \"256\u00f76=42, 4\"],\n  [\"991\u00f76=165, 1\", \"470\u00f74=117, 2\"],\n  [\"890\u00f78=111, 2\", \"473\u00f74=118, 1\"],\n  [\"491\u00f74=122, 3\", \"301\u00f77=43, 0\"],\n  [\"258\u00f74=64, 2\", \"795\u00f73=265, 0\"],\n  [\"334\u00f74=83, 2\", \"802\u00f74=200, 2\"],\n  [\"908\u00f75=181, 3\", \"810\u00f79=90, 0\"],\n  [\"652\u00f72=326, 0\", \"286\u00f77=40, 6\"],\n  [\"337\u00f72=168, 1\", \"715\u00f78=89, 3\"],\n  [\"528\u00f79=58, 6\", \"512\u00f73=170, 2\"],\n  [\"483\u00f78=60, 3\", \"345\u00f79=38, 3\"],\n  [\"545\u00f73=181, 2\", \"516\u00f72=258, 0\"],\n  [\"969\u00f76=161, 3\", \"218\u00f76=36, 2\"],\n  [\"231\u00f78=28, 7\", \"164\u00f77=23, 3\"],\n  [\"378\u00f72=189, 0\", \"226\u00f77=32, 2\"],\n  [\"986\u00f72=493, 0\", \"449\u00f76=74, 5\"],\n  [\"880\u00f78=110, 0\", \"314\u00f78=39, 2\"],\n  [\"256\u00f73=85, 1\", \"825\u00f75=165, 0\"],\n  [\"404\u00f77=57, 5\", \"236\u00f77=33, 5\"],\n  [\"578\u00f72=289, 0\", \"284\u00f79=31, 5\"],\n  [\"330\u00f78=41, 2\", \"509\u00f79=56, 5\"],\n  [\"779\u00f78=97, 3\", \"967\u00f73=322, 1\"],\n  [\"731\u00f78=91, 3\", \"627\u00f75=125, 2\"],\n  [\"298\u00f78=37, 2\", \"345\u00f74=86, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date and each \"three-digit \u00f7 one-digit\" practice answer\n# with the updated values, as described in the commit diff.\n# Each old value is unique in the document, so Find/Replace with\n# MatchCase and ReplaceAll is sufficient and leaves the surrounding\n# run formatting (font/size) untouched.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-02-13 Tuesday\", \"2024-02-14 Wednesday\"),\n    @(\"160\u00f74=40, 0\", \"979\u00f75=195, 4\"),\n    @(\"629\u00f78=78, 5\", \"256\u00f76=42, 4\"),\n    @(\"991\u00f76=165, 1\", \"470\u00f74=117, 2\"),\n    @(\"890\u00f78=111, 2\", \"473\u00f74=118, 1\"),\n    @(\"491\u00f74=122, 3\", \"301\u00f77=43, 0\"),\n    @(\"258\u00f74=64, 2\", \"795\u00f73=265, 0\"),\n    @(\"334\u00f74=83, 2\", \"802\u00f74=200, 2\"),\n    @(\"908\u00f75=181, 3\", \"810\u00f79=90, 0\"),\n    @(\"652\u00f72=326, 0\", \"286\u00f77=40, 6\"),\n    @(\"337\u00f72=168, 1\", \"715\u00f78=89, 3\"),\n    @(\"528\u00f79=58, 6\", \"512\u00f73=170, 2\"),\n    @(\"483\u00f78=60, 3\", \"345\u00f79=38, 3\"),\n    @(\"545\u00f73=181, 2\", \"516\u00f72=258, 0\"),\n    @(\"969\u00f76=161, 3\", \"218\u00f76=36, 2\"),\n    @(\"231\u00f78=28, 7\", \"164\u00f77=23, 3\"),\n    @(\"378\u00f72=189, 0\", \"226\u00f77=32, 2\"),\n    @(\"986\u00f72=493, 0\", \"449\u00f76=74, 5\"),\n    @(\"880\u00f78=110, 0\", \"314\u00f78=39, 2\"),\n    @(\"256\u00f73=85, 1\", \"825\u00f75=165, 0\"),\n    @(\"404\u00f77=57, 5\", \"236\u00f77=33, 5\"),\n    @(\"578\u00f72=289, 0\", \"284\u00f79=31, 5\"),\n    @(\"330\u00f78=41, 2\", \"509\u00f79=56, 5\"),\n    @(\"779\u00f78=97, 3\", \"967\u00f73=322, 1\"),\n    @(\"731\u00f78=91, 3\", \"627\u00f75=125, 2\"),\n    @(\"298\u00f78=37, 2\", \"345\u00f74=86, 1\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($old, $false, $true, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
